$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'29.004.20"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'  -1.82%  "
$ws.Range("E2").Style = "Normal"

$ws.Range("D3").Value = "'1.907.48"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'  -4.04%  "
$ws.Range("E3").Style = "Normal"

$ws.Range("D4").Value = "'1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'  -0.24%  "
$ws.Range("E4").Style = "Normal"

$ws.Range("D5").Value = "'324.59"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'  -1.48%  "
$ws.Range("E5").Style = "Normal"

$ws.Range("D6").Value = "'1.000"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'  -0.44%  "
$ws.Range("E6").Style = "Normal"

$ws.Range("D7").Value = "'0.4598"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'  -1.71%  "
$ws.Range("E7").Style = "Normal"

$ws.Range("D8").Value = "'0.3833"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'  -2.84%  "
$ws.Range("E8").Style = "Normal"

$ws.Range("D9").Value = "'0.07752"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'  -2.72%  "
$ws.Range("E9").Style = "Normal"

$ws.Range("D10").Value = "'0.9860"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'  -1.71%  "
$ws.Range("E10").Style = "Normal"

$ws.Range("D11").Value = "'22.10"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'  -3.09%  "
$ws.Range("E11").Style = "Normal"

$ws.Range("D12").Value = "'1.910.60"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'  -5.39%  "
$ws.Range("E12").Style = "Normal"

$ws.Range("D13").Value = "'7.004"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'  -3.54%  "
$ws.Range("E13").Style = "Normal"

$ws.Range("D14").Value = "'5.706"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'  -2.97%  "
$ws.Range("E14").Style = "Normal"

$ws.Range("D15").Value = "'0.07060"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'  -1.48%  "
$ws.Range("E15").Style = "Normal"

$ws.Range("E16").Value = "'  -0.48%  "
$ws.Range("E16").Style = "Normal"

$ws.Range("D17").Value = "'84.12"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'  -5.35%  "
$ws.Range("E17").Style = "Normal"

$ws.Range("D18").Value = "'0.000009569"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'  -4.05%  "
$ws.Range("E18").Style = "Normal"

$ws.Range("D19").Value = "'16.74"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'  -3.59%  "
$ws.Range("E19").Style = "Normal"

$ws.Range("E20").Value = "'  -0.45%  "
$ws.Range("E20").Style = "Normal"

$ws.Range("D21").Value = "'29.012.97"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "'  -2.07%  "
$ws.Range("E21").Style = "Normal"

$ws.Range("D22").Value = "'5.330"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'  -3.96%  "
$ws.Range("E22").Style = "Normal"

$ws.Range("D23").Value = "'10.97"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'  -2.81%  "
$ws.Range("E23").Style = "Normal"

$ws.Range("D24").Value = "'2.135.32"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "'  -5.16%  "
$ws.Range("E24").Style = "Normal"

$ws.Range("D25").Value = "'2.074"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'  -2.29%  "
$ws.Range("E25").Style = "Normal"

$ws.Range("D26").Value = "'156.32"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'  -1.29%  "
$ws.Range("E26").Style = "Normal"

$ws.Range("D27").Value = "'19.22"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'  -2.30%  "
$ws.Range("E27").Style = "Normal"

$ws.Range("D28").Value = "'5.630"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "'  -5.81%  "
$ws.Range("E28").Style = "Normal"

$ws.Range("D29").Value = "'117.93"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "'  -2.06%  "
$ws.Range("E29").Style = "Normal"

$ws.Range("D30").Value = "'1.832"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "'  -6.78%  "
$ws.Range("E30").Style = "Normal"

$ws.Range("D31").Value = "'0.09265"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "'  -1.91%  "
$ws.Range("E31").Style = "Normal"

$ws.Range("D32").Value = "'0.8634"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "'  -3.35%  "
$ws.Range("E32").Style = "Normal"

$ws.Range("D33").Value = "'5.120"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "'  -3.16%  "
$ws.Range("E33").Style = "Normal"

$ws.Range("D34").Value = "'1.254"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "'  -6.92%  "
$ws.Range("E34").Style = "Normal"

$ws.Range("D35").Value = "'3.015"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "'  -5.52%  "
$ws.Range("E35").Style = "Normal"

$ws.Range("D36").Value = "'0.05737"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "'  -1.76%  "
$ws.Range("E36").Style = "Normal"

$ws.Range("E37").Value = "'  -2.11%  "
$ws.Range("E37").Style = "Normal"

$ws.Range("E38").Value = "'  -0.72%  "
$ws.Range("E38").Style = "Normal"

$ws.Range("D39").Value = "'0.02052"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "'  -3.82%  "
$ws.Range("E39").Style = "Normal"

$ws.Range("D40").Value = "'7.491"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'  -5.46%  "
$ws.Range("E40").Style = "Normal"

$ws.Range("D41").Value = "'0.5541"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'  -3.86%  "
$ws.Range("E41").Style = "Normal"

$ws.Range("D42").Value = "'0.1760"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'  -3.55%  "
$ws.Range("E42").Style = "Normal"

$ws.Range("D43").Value = "'9.302"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "'  -5.35%  "
$ws.Range("E43").Style = "Normal"

$ws.Range("D44").Value = "'2.724"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'  +3.02%  "
$ws.Range("E44").Style = "Normal"

$ws.Range("D45").Value = "'0.000002728"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "'  -12.50%  "
$ws.Range("E45").Style = "Normal"

$ws.Range("D46").Value = "'0.5222"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "'  -3.03%  "
$ws.Range("E46").Style = "Normal"

$ws.Range("D47").Value = "'11.35"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "'  -5.92%  "
$ws.Range("E47").Style = "Normal"

$ws.Range("D48").Value = "'2.115"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "'  -2.73%  "
$ws.Range("E48").Style = "Normal"

$ws.Range("D49").Value = "'0.06832"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "'  -1.91%  "
$ws.Range("E49").Style = "Normal"

$ws.Range("D50").Value = "'112.13"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "'  -2.29%  "
$ws.Range("E50").Style = "Normal"

$ws.Range("E51").Value = "'  -4.62%  "
$ws.Range("E51").Style = "Normal"
